$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 7.545
    "B6"  = 5.551
    "B7"  = 5.396000000000001
    "E7"  = 16.128
    "E12" = 17.638
    "E15" = 16.113
    "B16" = 5.429
    "B20" = 6.891000000000001
    "E20" = 16.291
    "E21" = 16.523
    "E22" = 16.54
    "E23" = 16.527
    "B28" = 6.171000000000001
    "B29" = 5.306
    "E29" = 17.246
    "B32" = 6.453
    "E34" = 16.747
    "B40" = 9.343999999999999
    "E42" = 16.539
    "E43" = 16.949
    "E44" = 16.388
    "E45" = 16.774
    "B46" = 6.382000000000001
    "E46" = 16.659
    "E50" = 16.535
    "B51" = 5.884
    "E51" = 16.957
    "B52" = 5.369
    "B57" = 5.189
    "B59" = 4.607999999999999
    "B62" = 5.261
    "B66" = 5.013
    "E66" = 17.384
    "E67" = 17.571
    "B73" = 6.503
    "B74" = 9.204000000000001
    "E79" = 16.957
    "E84" = 16.606
    "B92" = 5.198
    "E92" = 17.815
    "E97" = 16.872
    "B100" = 5.849
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
